$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the column header suffixes: "_old" -> "_FV2404", "_new" -> "_FV2410"
#    (the header row, A1:U1, are the only cells carrying these suffixes).
$ws.Cells.Replace("_old", "_FV2404", 2)
$ws.Cells.Replace("_new", "_FV2410", 2)

# 2. Turn the data range into an Excel Table ("Table1") with an AutoFilter.
$tableRange = $ws.Range("A1:U66")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# 3. Freeze the header row (split/freeze at row 2, i.e. above row 2).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
